# edit.ps1 - apply the "fixed report and presentation" commit to the deck.
#
# Summary of the change being reproduced:
#   * Slide 3 & slide 4 title placeholders: drop the 90% autofit font-scale
#     (<a:normAutofit fontScale="90000"/> -> <a:normAutofit/>).
#   * Slide 4's content placeholder is moved to the top-left position that
#     slides 5/6 will also use.
#   * Two new "Title and Content" slides are appended (slide 5 "Klašu
#     diagramma", slide 6 "Galveno metožu algoritmu apraksts"), cloned from
#     slide 4 so they inherit identical placeholder styling.
#   * The theme's major/minor latin font is switched to Cambria.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: PowerPoint's Shape.Left/Top/Width/Height (and friends) are in
# points, while OOXML stores EMU (1 pt = 12700 EMU). Converting EMU/12700.0
# and writing it back loses the last EMU to float rounding inside the COM
# bridge, so nudge by half an EMU (expressed in points) before handing the
# value to the property setter - that is enough to land back on the exact
# integer EMU value on write-back without affecting any real placement.
function EMU([double]$emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

# ---------------------------------------------------------------------------
# 1. Duplicate slide 4 twice to become the new slides 5 and 6. Duplicating
#    from the growing tail (instead of always from slide 4) keeps the new
#    slide IDs in ascending order, matching how PowerPoint numbers them.
$null = $p.Slides.Item(4).Duplicate()
$null = $p.Slides.Item(5).Duplicate()

$slide3 = $p.Slides.Item(3)
$slide4 = $p.Slides.Item(4)
$slide5 = $p.Slides.Item(5)
$slide6 = $p.Slides.Item(6)

# ---------------------------------------------------------------------------
# 2. Slide 3 ("Sistēmas funkcionālās prasības"): title no longer shrinks to
#    90% - switching AutoSize through ppAutoSizeTextToFitShape recomputes
#    the <a:normAutofit/> without a leftover fontScale.
$slide3.Shapes.Item(1).TextFrame.AutoSize = 2

# ---------------------------------------------------------------------------
# 3. Slide 4 ("Sistēmas nefunkcionālās prasības"): same title fix, plus the
#    content placeholder moves up to the new top-left spot.
$slide4.Shapes.Item(1).TextFrame.AutoSize = 2
$slide4.Shapes.Item(2).Left = EMU 372533
$slide4.Shapes.Item(2).Top = EMU 1394989

# ---------------------------------------------------------------------------
# 4. Slide 5 (new): "Klašu diagramma" - content placeholder in the new spot,
#    title keeps the de-scaled autofit (same fix as slides 3/4).
$slide5.Shapes.Item(1).TextFrame.TextRange.Text = "Klašu diagramma"
$slide5.Shapes.Item(1).TextFrame.AutoSize = 2
$slide5.Shapes.Item(2).Left = EMU 372533
$slide5.Shapes.Item(2).Top = EMU 1394989

# ---------------------------------------------------------------------------
# 5. Slide 6 (new): "Galveno metožu algoritmu apraksts" - content placeholder
#    moves too, but the title keeps its inherited 90% fontScale autofit.
$slide6.Shapes.Item(1).TextFrame.TextRange.Text = "Galveno metožu algoritmu apraksts"
$slide6.Shapes.Item(2).Left = EMU 372533
$slide6.Shapes.Item(2).Top = EMU 1394989

# ---------------------------------------------------------------------------
# 6. Theme: switch the major/minor latin typeface to Cambria.
$fontScheme = $p.SlideMaster.Theme.ThemeFontScheme
$fontScheme.MajorFont.Latin = "Cambria"
$fontScheme.MinorFont.Latin = "Cambria"
